$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Critical_Points_Follow_Up")

$ws.Range("B3").Value = "1"
$ws.Range("B4").Value = "2"
$ws.Range("B5").Value = "3"

$ws.Range("C3").Value = "A"
$ws.Range("C4").Value = "B"
$ws.Range("C5").Value = "C"
$ws.Range("C6").Value = "D"

$ws.Range("D3").Value = "Open"
$ws.Range("D4").Value = "Closed"
$ws.Range("D5").Value = "Closed"
$ws.Range("D6").Value = "In Progress"

$ws.Range("D7").Select()
